# Implement suggestion as per coach request
#
# category sheet: add a "Type" header, rename a couple of categories, and
# add a new "Feria" (Expense) category.
#
# "Expense and incomes" sheet: bump the university fee amount, replace the
# old "Servicios Profesionales" income entry with a new "Feria" expense
# entry, and re-add the professional-services income entry (now lower-case
# detail / underscored category) as a new row with an updated amount.

$wb = $excel.ActiveWorkbook

$category = $wb.Worksheets.Item("category")
$expenses = $wb.Worksheets.Item("Expense and incomes")

# --- category sheet -------------------------------------------------
$category.Range("B1").Value = "Type"

$category.Range("A2").Value = "Casa"

$category.Range("A4").Value = "Servicios_Profesionales"

$category.Range("A5").Value = "Feria"
$category.Range("B5").Value = "Expense"

# --- Expense and incomes sheet --------------------------------------
# Row 2: tuition amount increased. Keep it stored as text (matches the
# original inline-string typing) via a quote-prefixed literal, then drop
# the quote-prefix formatting Excel applies so the cell keeps the sheet's
# default (unstyled) look.
$expenses.Range("D2").Value = "'47500"
$expenses.Range("D2").ClearFormats()

# Row 3: was the "Servicios Profesionales" income row, now becomes the
# "Feria" expense row with its own detail/amount.
$expenses.Range("A3").Value = "Tomate,Zanahoria,Culantro"
$expenses.Range("B3").Value = "Feria"
$expenses.Range("C3").Value = "Expense"
$expenses.Range("D3").Value = "'1500"
$expenses.Range("D3").ClearFormats()

# Row 4 (new): professional-services income, re-added with the new
# underscored category name and the updated amount.
$expenses.Range("A4").Value = "analisis de un sistema informatico"
$expenses.Range("B4").Value = "Servicios_Profesionales"
$expenses.Range("C4").Value = "Income"
$expenses.Range("D4").Value = "'1500"
$expenses.Range("D4").ClearFormats()

# Make "Expense and incomes" the active sheet/selection, matching the
# coach-requested view state.
$expenses.Activate()
$expenses.Range("A3").Select() | Out-Null

Write-Host "Applied coach-requested category/expense updates."
